$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll the sheet view back to the top (topLeftCell C5 -> C1)
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 3

# New query function made for row 7 (EXONS_POSITIONS) -> mark Status as Done
$ws.Range("F7").Value = "Done"
$ws.Range("F7").Style = $ws.Range("F6").Style

